$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column L (2020) mirroring column K's styles/values for rows 4-12
$ws.Range("L4").Value = 2020
$ws.Range("L5").Value = 5.6
$ws.Range("L6").Value = 0.8
$ws.Range("L7").Value = 1.9
$ws.Range("L8").Value = 0.7
$ws.Range("L9").Value = 0.7
$ws.Range("L10").Value = 0.9
$ws.Range("L11").Value = 0.3
$ws.Range("L12").Value = 0.2

# Copy styles from column K to column L for the same rows
$ws.Range("K4:K12").Copy()
$ws.Range("L4:L12").PasteSpecial(-4122)  # xlPasteFormats

# Update selection to match target state
$ws.Range("N5").Select()
